$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet")
$ws2 = $wb.Worksheets.Item("Evaluation Warning")

# Update the DevExpress evaluation-warning version string
$ws2.Range("A6").Value = "or purchase a new license (devexpress.com/BUY) to continue use of DevExpress product libraries (v25.2.3.0)."

# Update data rows on the main "Sheet"
$ws1.Range("A2").Value = 482
$ws1.Range("I2").Value = 252
$ws1.Range("J2").Value = 46025.683620825424

$ws1.Range("A3").Value = 483
$ws1.Range("I3").Value = 253
$ws1.Range("J3").Value = 46025.683620825424
